$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.11"
$ws.Range("G2").Value = "'11"
$ws.Range("D3").Value = "'22.39"
$ws.Range("G3").Value = "'11"
$ws.Range("D4").Value = "'5.511"
$ws.Range("G4").Value = "'11"
$ws.Range("D5").Value = "'0.05622"
$ws.Range("G5").Value = "'11"
$ws.Range("D6").Value = "'6.467"
$ws.Range("G6").Value = "'11"
$ws.Range("G7").Value = "'11"
$ws.Range("D8").Value = "'1.048"
$ws.Range("G8").Value = "'11"
$ws.Range("D9").Value = "'0.1425"
$ws.Range("G9").Value = "'11"
$ws.Range("D10").Value = "'0.07279"
$ws.Range("G10").Value = "'11"
$ws.Range("D11").Value = "'0.03178"
$ws.Range("G11").Value = "'11"
$ws.Range("D12").Value = "'0.02976"
$ws.Range("G12").Value = "'11"
$ws.Range("D13").Value = "'0.09260"
$ws.Range("G13").Value = "'11"
$ws.Range("D14").Value = "'0.001670"
$ws.Range("G14").Value = "'11"
$ws.Range("D15").Value = "'3.203"
$ws.Range("G15").Value = "'11"
$ws.Range("D16").Value = "'0.04708"
$ws.Range("G16").Value = "'11"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006278"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("G17").Value = "'11"
$ws.Range("B18").Value = "BitKan"
$ws.Range("C18").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D18").Value = "'0.001051"
$ws.Range("E18").Value = "17BitKanKAN"
$ws.Range("G18").Value = "'11"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Value = "'0.003819"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("G19").Value = "'11"
$ws.Range("B20").Value = "NitroEx"
$ws.Range("C20").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D20").Value = "'0.0001503"
$ws.Range("E20").Value = "19NitroExNTX"
$ws.Range("G20").Value = "'11"
$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").Value = "'0.0003306"
$ws.Range("E21").Value = "20UpBotsUBXT"
$ws.Range("G21").Value = "'11"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").Value = "'3.976"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("G22").Value = "'11"
$ws.Range("B23").Value = "GateToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D23").Value = "'3.391"
$ws.Range("E23").Value = "22GateTokenGT"
$ws.Range("G23").Value = "'11"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.120"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("G24").Value = "'11"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Value = "'0.0005982"
$ws.Range("E25").Value = "24OneONE"
$ws.Range("G25").Value = "'11"
$ws.Range("G26").Value = "'11"
$ws.Range("E27").Value = "26ProBitTokenPROB"
$ws.Range("G27").Value = "'11"
$ws.Range("G28").Value = "'11"
$ws.Range("G29").Value = "'11"
$ws.Range("G30").Value = "'11"
$ws.Range("G31").Value = "'11"
$ws.Range("G32").Value = "'11"
$ws.Range("G33").Value = "'11"
$ws.Range("G34").Value = "'11"
$ws.Range("G35").Value = "'11"
$ws.Range("G36").Value = "'11"
$ws.Range("G37").Value = "'11"
$ws.Range("G38").Value = "'11"
$ws.Range("G39").Value = "'11"
$ws.Range("G40").Value = "'11"
$ws.Range("D41").Value = "'0.006893"
$ws.Range("G41").Value = "'11"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003507"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = "'11"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1042"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = "'11"
$ws.Range("D44").Value = "'0.01022"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("G44").Value = "'11"
$ws.Range("D45").Value = "'0.00005647"
$ws.Range("G45").Value = "'11"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("G46").Value = "'11"
$ws.Range("D47").Value = "'0.6813"
$ws.Range("G47").Value = "'11"
$ws.Range("D48").Value = "'0.02581"
$ws.Range("G48").Value = "'11"
$ws.Range("G49").Value = "'11"
$ws.Range("G50").Value = "'11"
$ws.Range("G51").Value = "'11"
